# edit.ps1 - applies the funding-slide reshuffle + new References link
# described in the commit diff.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 7 - "Who can fund and help us achieve these goals?"
# Content Placeholder 2: insert two new bullet paragraphs ("Tesla" and
# "INEC - organize the elections in Nigeria") right after "CISCO" and
# before the crowdfunding paragraph.
# ---------------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
$fundingShape = $slide7.Shapes.Item(4)
$fundingRange = $fundingShape.TextFrame.TextRange

$fundingParas = $fundingRange.Paragraphs()
$ciscoPara = $fundingParas.Paragraphs(3, 1)

$cr = [char]13
$dash = [char]0x2013
$ciscoPara.InsertAfter($cr + "Tesla" + $cr + "INEC " + $dash + " organize the elections in Nigeria")

# ---------------------------------------------------------------------
# Slide 9 - "References"
# Title shape: rebuild the paragraph so the stray trailing endParaRPr
# (left over from the old text) is gone, matching a freshly retyped run.
# ---------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$titleShape = $slide9.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Delete()
$titleShape.TextFrame.TextRange.InsertAfter("References")

# Content Placeholder 2: add a third link, "www.google.com", after the
# CNN link, with its own hyperlink (new relationship rId4).
$linksShape = $slide9.Shapes.Item(3)
$linksRange = $linksShape.TextFrame.TextRange
$linksParas = $linksRange.Paragraphs()
$cnnPara = $linksParas.Paragraphs(2, 1)
$cnnPara.InsertAfter([char]13 + "www.google.com")

$linksParas2 = $linksShape.TextFrame.TextRange.Paragraphs()
$googlePara = $linksParas2.Paragraphs(3, 1)
$googleAction = $googlePara.ActionSettings.Item(1)
$googleAction.Hyperlink.Address = "http://www.google.com"
